# Append four new sales rows (2023-10-31 entries for Oral B and Waw Detergent)
# to the bottom of the existing sales table on the active sheet, extending the
# dimension from A1:F13 to A1:F17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows, in column order: PC, DATE, PRD, SP, QTY, AMT
$newRows = @(
    @("121210", "2023-10-31", "Oral B",         "800",  "3", "2400"),
    @("121210", "2023-10-31", "Oral B",         "800",  "6", "4800"),
    @("121213", "2023-10-31", "Waw Detergent",  "1000", "4", "4000"),
    @("121213", "2023-10-31", "Waw Detergent",  "1000", "5", "5000")
)

$startRow = 14
$endRow = $startRow + $newRows.Count - 1

# Every existing cell in the sheet is stored as text (PC codes, dates, qty,
# amounts are all plain strings), even though they look numeric. Mark the
# target range as Text *before* writing so Excel doesn't auto-convert values
# like "121210" or "2023-10-31" into a number/date.
$targetRange = $ws.Range("A" + $startRow + ":F" + $endRow)
$targetRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $newRows[$i][$col - 1]
    }
}

# Restore the default (General) style on the new cells so they match the
# rest of the sheet, which carries no explicit cell style.
$targetRange.Style = $ws.Range("A2:F2").Style
